$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.190.65'
$ws.Range("E2").Value = '  +4.38%  '
$ws.Range("D3").Value = '1.687.12'
$ws.Range("E3").Value = '  +3.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.44'
$ws.Range("E5").Value = '  +2.55%  '
$ws.Range("E6").Value = '  +2.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.49'
$ws.Range("E8").Value = '  +2.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.267'
$ws.Range("E9").Value = '  +3.24%  '
$ws.Range("E10").Value = '  +5.04%  '
$ws.Range("D12").Value = '1.925.91'
$ws.Range("E12").Value = '  +3.31%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.40'
$ws.Range("E13").Value = '  +10.57%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.680.58'
$ws.Range("E14").Value = '  +3.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.611'
$ws.Range("E15").Value = '  +7.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.15'
$ws.Range("E16").Value = '  +8.26%  '
$ws.Range("D17").Value = '31.124.55'
$ws.Range("E17").Value = '  +4.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.79'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '249.73'
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("D20").Value = '0.0₃0723'
$ws.Range("E20").Value = '  +2.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.29'
$ws.Range("E22").Value = '  +3.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.06'
$ws.Range("E23").Value = '  +2.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.73'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.99'
$ws.Range("E26").Value = '  +3.20%  '
$ws.Range("E27").Value = '  +2.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.73'
$ws.Range("E28").Value = '  +1.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0499'
$ws.Range("E30").Value = '  +1.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.53'
$ws.Range("E31").Value = '  +4.36%  '
$ws.Range("E32").Value = '  +3.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.35'
$ws.Range("E33").Value = '  +5.13%  '
$ws.Range("D34").Value = '1.520.36'
$ws.Range("E34").Value = '  +6.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '84.67'
$ws.Range("E36").Value = '  +12.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.03'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("E38").Value = '  +10.35%  '
$ws.Range("E39").Value = '  +5.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.69'
$ws.Range("E40").Value = '  -2.67%  '
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("E42").Value = '  +3.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.845'
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0505'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.59'
$ws.Range("E47").Value = '  +8.29%  '
$ws.Range("E48").Value = '  +4.82%  '
$ws.Range("D49").Value = '1.817.09'
$ws.Range("E49").Value = '  +2.57%  '
$ws.Range("E50").Value = '  +8.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '93.69'
$ws.Range("E51").Value = '  +1.79%  '
